$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Razon social" entries: commas used as separators changed to periods ---
$ws.Range("E65").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E76").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E89").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E96").Value = "MONROY. AGUSTIN ALEJANDRO"
$ws.Range("E174").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E178").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E180").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E197").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix "Importe" column: Spanish-formatted numeric text (1.234,56) -> plain decimal text (1234.56) ---
# Leading apostrophe forces Excel to keep the value as literal text (quote-prefixed),
# matching how these amounts are already stored as text strings, not numbers.
$importe = @{
    2 = "2630.00"
    3 = "5300.00"
    4 = "38720.00"
    5 = "38720.00"
    6 = "17325.00"
    7 = "170.00"
    8 = "13407.03"
    9 = "1680.00"
    10 = "5759.98"
    11 = "246.11"
    12 = "599.20"
    13 = "13967.31"
    14 = "154627.52"
    15 = "284336.32"
    16 = "41789.20"
    17 = "15960.00"
    18 = "98062.43"
    19 = "2475.00"
    20 = "9268.93"
    21 = "4920.00"
    22 = "33347.00"
    23 = "45874.00"
    24 = "20213.69"
    25 = "28788.74"
    26 = "24336.70"
    27 = "950.00"
    28 = "5840.00"
    29 = "318.46"
    30 = "470.00"
    31 = "17.40"
    32 = "33956.04"
    33 = "157800.00"
    34 = "71.64"
    35 = "21873.21"
    36 = "4751.92"
    37 = "2962.06"
    38 = "567.60"
    39 = "1760.00"
    40 = "20435.40"
    41 = "26681.00"
    42 = "4126.97"
    43 = "26427.90"
    44 = "5998.00"
    45 = "3600.00"
    46 = "60477.82"
    47 = "146.00"
    48 = "419.00"
    49 = "16222.79"
    50 = "117.12"
    51 = "2961.20"
    52 = "287908.32"
    53 = "1342.81"
    54 = "5450.00"
    55 = "1656.99"
    56 = "110.10"
    57 = "890.40"
    58 = "500.00"
    59 = "298.35"
    60 = "5731.51"
    61 = "388.36"
    62 = "19314.00"
    63 = "89176.21"
    64 = "444.00"
    65 = "760.00"
    66 = "10215.92"
    67 = "410.37"
    68 = "2535.00"
    69 = "352.05"
    70 = "7297.00"
    71 = "12351.00"
    72 = "690.00"
    73 = "7654.40"
    74 = "8668.00"
    75 = "2200.00"
    76 = "1240.00"
    77 = "1449.66"
    78 = "28347.92"
    79 = "140.00"
    80 = "5835.00"
    81 = "51000.00"
    82 = "2400.00"
    83 = "2560.00"
    84 = "76290.00"
    85 = "22000.00"
    86 = "429.00"
    87 = "15267.00"
    88 = "2520.00"
    89 = "2550.00"
    90 = "390.60"
    91 = "792.00"
    92 = "287.50"
    93 = "1.20"
    94 = "834.78"
    95 = "415230.93"
    96 = "20300.00"
    97 = "53261.49"
    98 = "51015.00"
    99 = "3450.00"
    100 = "3050.00"
    101 = "48.80"
    102 = "550.00"
    103 = "44881.65"
    104 = "17.50"
    105 = "14369.00"
    106 = "21.00"
    107 = "47.50"
    108 = "3540.00"
    109 = "1100.00"
    110 = "4279.00"
    111 = "43225.00"
    112 = "1146.00"
    113 = "362.00"
    114 = "18400.00"
    115 = "537.00"
    116 = "1050.00"
    117 = "2400.00"
    118 = "765.00"
    119 = "125.00"
    120 = "32822.17"
    121 = "2336.00"
    122 = "2792.92"
    123 = "2155.82"
    124 = "943.00"
    125 = "3450.00"
    126 = "1407.80"
    127 = "90.00"
    128 = "31470.00"
    129 = "23608.25"
    130 = "1350.00"
    131 = "200430.00"
    132 = "180000.00"
    133 = "380.00"
    134 = "5584.00"
    135 = "3211.80"
    136 = "640.00"
    137 = "12430.00"
    138 = "17300.00"
    139 = "5000.00"
    140 = "15000.00"
    141 = "4000.00"
    142 = "1500.00"
    143 = "804.80"
    144 = "4852.15"
    145 = "142.80"
    146 = "2222.40"
    147 = "261.21"
    148 = "205500.00"
    149 = "29377.88"
    150 = "8000.00"
    151 = "8200.00"
    152 = "12000.00"
    153 = "4000.00"
    154 = "10000.00"
    155 = "2500.00"
    156 = "2762.50"
    157 = "2556.00"
    158 = "3000.00"
    159 = "2000.00"
    160 = "2000.00"
    161 = "1500.00"
    162 = "18426.50"
    163 = "7500.00"
    164 = "4000.00"
    165 = "2500.00"
    166 = "3800.00"
    167 = "4500.00"
    168 = "6770.00"
    169 = "9000.00"
    170 = "2500.00"
    171 = "1500.00"
    172 = "18320.00"
    173 = "1210.00"
    174 = "1000.00"
    175 = "51630.00"
    176 = "1200.00"
    177 = "1896.80"
    178 = "1200.00"
    179 = "399.67"
    180 = "13247.00"
    181 = "28805.00"
    182 = "7050.00"
    183 = "1380.00"
    184 = "53.14"
    185 = "250.00"
    186 = "734.70"
    187 = "13500.00"
    188 = "167.00"
    189 = "2809.00"
    190 = "113827.00"
    191 = "127.20"
    192 = "71.60"
    193 = "2390.00"
    194 = "119164.00"
    195 = "5746.32"
    196 = "15560.46"
    197 = "4310.00"
    198 = "5200.00"
    199 = "4565.99"
    200 = "2154.10"
    201 = "101144.77"
    202 = "21060.00"
    203 = "1400.00"
    204 = "212628.00"
    205 = "14092.00"
    206 = "47801.00"
    207 = "780.00"
    208 = "372.40"
    209 = "513795.21"
    210 = "10700.00"
    211 = "1475265.98"
    212 = "320236.00"
    213 = "6800.00"
    214 = "4500.00"
    215 = "244850.00"
    216 = "203500.00"
    217 = "262350.00"
    218 = "194336.00"
    219 = "32500.00"
    220 = "27000.00"
    221 = "227390.00"
    222 = "318372.00"
    223 = "220000.00"
    224 = "130000.00"
    225 = "167700.00"
    226 = "92780.00"
    227 = "110000.00"
    228 = "1087502.52"
    229 = "379855.81"
    230 = "33250.00"
    231 = "10286.40"
    232 = "17000.00"
    233 = "4500.00"
    234 = "17516.56"
    235 = "1100.00"
    236 = "92980.00"
    237 = "840.00"
    238 = "800.00"
}
foreach ($row in $importe.Keys) {
    $ws.Cells.Item($row, 8).Value = "'" + $importe[$row]
}
